$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = "FC1"
$ws.Range("G4").Value = "FC2"
$ws.Range("H4").Value = "FC3"

$ws.Range("F5").Value = 181
$ws.Range("G5").Value = 39
$ws.Range("H5").Value = 4

$ws.Range("H6").Select()
